$wb = $excel.ActiveWorkbook

$oldGuid = "4cc044c5-93a8-4f3c-9be7-ec7391096bff"
$newGuid = "a8427597-0e66-47f0-9acf-a73325e27a65"
$newHash = "c1145e64bb3279f92b7cb704bcd8bd0fdd43fce3"

$newFileName   = "$newGuid.md"
$newPathName   = "e2e\$newGuid.md"
$newHoDate     = "2016-08-26 15:13:37"
$newZhHandoff  = "$newGuid.$newHash.zh-cn.xlf"
$newDeHandoff  = "$newGuid.$newHash.de-de.xlf"
$newZhHandoffDate = "2016-08-26 15:13:33"
$neverDate     = "0001-01-01 00:00:00"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathName
$wsOverview.Range("G2").Value = $newHoDate

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newFileName
$wsZh.Range("G2").Value = $newZhHandoff
$wsZh.Range("H2").Value = $newZhHandoffDate
$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $neverDate

$wsZh.Columns.Item(9).AutoFit() | Out-Null
$wsZh.Columns.Item(10).AutoFit() | Out-Null
$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newFileName
$wsDe.Range("G2").Value = $newDeHandoff
$wsDe.Range("H2").Value = $newHoDate
$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $neverDate

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
